$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Retrained model data: update existing rows 2-59 (timestamps shift from
# 2025-... day to the next day, plus new Import/Export forecasts for ELnet/NRG)
$ws.Cells.Item(2, 1).Value = 45982
$ws.Cells.Item(2, 2).Value = 2.434
$ws.Cells.Item(2, 3).Value = 6.682
$ws.Cells.Item(3, 1).Value = 45982.01041666666
$ws.Cells.Item(3, 2).Value = 13.159
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(4, 1).Value = 45982.02083333334
$ws.Cells.Item(4, 2).Value = 14.898
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(5, 1).Value = 45982.03125
$ws.Cells.Item(5, 2).Value = 25.861
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(6, 1).Value = 45982.04166666666
$ws.Cells.Item(6, 2).Value = 10.729
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(7, 1).Value = 45982.05208333334
$ws.Cells.Item(7, 2).Value = 5.204
$ws.Cells.Item(7, 3).Value = 1.382
$ws.Cells.Item(8, 1).Value = 45982.0625
$ws.Cells.Item(8, 2).Value = 21.775
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(9, 1).Value = 45982.07291666666
$ws.Cells.Item(9, 2).Value = 21.957
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(10, 1).Value = 45982.08333333334
$ws.Cells.Item(10, 2).Value = 15.17
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(11, 1).Value = 45982.09375
$ws.Cells.Item(11, 2).Value = 21.101
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(12, 1).Value = 45982.10416666666
$ws.Cells.Item(12, 2).Value = 15.088
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(13, 1).Value = 45982.11458333334
$ws.Cells.Item(13, 2).Value = 15.57
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 1).Value = 45982.125
$ws.Cells.Item(14, 2).Value = 7.731
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(15, 1).Value = 45982.13541666666
$ws.Cells.Item(15, 2).Value = 19.647
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(16, 1).Value = 45982.14583333334
$ws.Cells.Item(16, 2).Value = 33.607
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(17, 1).Value = 45982.15625
$ws.Cells.Item(17, 2).Value = 42.088
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(18, 1).Value = 45982.16666666666
$ws.Cells.Item(18, 2).Value = 17.39
$ws.Cells.Item(18, 3).Value = 0.016
$ws.Cells.Item(19, 1).Value = 45982.17708333334
$ws.Cells.Item(19, 2).Value = 4.052
$ws.Cells.Item(19, 3).Value = 1.25
$ws.Cells.Item(20, 1).Value = 45982.1875
$ws.Cells.Item(20, 2).Value = 1.613
$ws.Cells.Item(20, 3).Value = 1.91
$ws.Cells.Item(21, 1).Value = 45982.19791666666
$ws.Cells.Item(21, 2).Value = 10.796
$ws.Cells.Item(21, 3).Value = 0.063
$ws.Cells.Item(22, 1).Value = 45982.20833333334
$ws.Cells.Item(22, 2).Value = 21.301
$ws.Cells.Item(22, 3).Value = 0.483
$ws.Cells.Item(23, 1).Value = 45982.21875
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 10.072
$ws.Cells.Item(24, 1).Value = 45982.22916666666
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(24, 3).Value = 8.683999999999999
$ws.Cells.Item(25, 1).Value = 45982.23958333334
$ws.Cells.Item(25, 2).Value = 0.779
$ws.Cells.Item(25, 3).Value = 1.774
$ws.Cells.Item(26, 1).Value = 45982.25
$ws.Cells.Item(26, 2).Value = 0.056
$ws.Cells.Item(26, 3).Value = 12.234
$ws.Cells.Item(27, 1).Value = 45982.26041666666
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 25.384
$ws.Cells.Item(28, 1).Value = 45982.27083333334
$ws.Cells.Item(28, 2).Value = 0.003
$ws.Cells.Item(28, 3).Value = 7.288
$ws.Cells.Item(29, 1).Value = 45982.28125
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(29, 3).Value = 14.045
$ws.Cells.Item(30, 1).Value = 45982.29166666666
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 22.737
$ws.Cells.Item(31, 1).Value = 45982.30208333334
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 25.189
$ws.Cells.Item(32, 1).Value = 45982.3125
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(32, 3).Value = 29.633
$ws.Cells.Item(33, 1).Value = 45982.32291666666
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 35.244
$ws.Cells.Item(34, 1).Value = 45982.33333333334
$ws.Cells.Item(34, 2).Value = 0.132
$ws.Cells.Item(34, 3).Value = 11.329
$ws.Cells.Item(35, 1).Value = 45982.34375
$ws.Cells.Item(35, 2).Value = 0
$ws.Cells.Item(35, 3).Value = 5.535
$ws.Cells.Item(36, 1).Value = 45982.35416666666
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(36, 3).Value = 13.592
$ws.Cells.Item(37, 1).Value = 45982.36458333334
$ws.Cells.Item(37, 2).Value = 0.623
$ws.Cells.Item(37, 3).Value = 5.076
$ws.Cells.Item(38, 1).Value = 45982.375
$ws.Cells.Item(38, 2).Value = 9.321
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(39, 1).Value = 45982.38541666666
$ws.Cells.Item(39, 2).Value = 0.5590000000000001
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(40, 1).Value = 45982.39583333334
$ws.Cells.Item(40, 2).Value = 0.116
$ws.Cells.Item(40, 3).Value = 18.794
$ws.Cells.Item(41, 1).Value = 45982.40625
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(41, 3).Value = 15.822
$ws.Cells.Item(42, 1).Value = 45982.41666666666
$ws.Cells.Item(42, 2).Value = 16.546
$ws.Cells.Item(42, 3).Value = 0.634
$ws.Cells.Item(43, 1).Value = 45982.42708333334
$ws.Cells.Item(43, 2).Value = 28.682
$ws.Cells.Item(43, 3).Value = 0
$ws.Cells.Item(44, 1).Value = 45982.4375
$ws.Cells.Item(44, 2).Value = 1.393
$ws.Cells.Item(44, 3).Value = 3.243
$ws.Cells.Item(45, 1).Value = 45982.44791666666
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(45, 3).Value = 35.148
$ws.Cells.Item(46, 1).Value = 45982.45833333334
$ws.Cells.Item(46, 2).Value = 0
$ws.Cells.Item(46, 3).Value = 27.383
$ws.Cells.Item(47, 1).Value = 45982.46875
$ws.Cells.Item(47, 2).Value = 0
$ws.Cells.Item(47, 3).Value = 40.888
$ws.Cells.Item(48, 1).Value = 45982.47916666666
$ws.Cells.Item(48, 2).Value = 0
$ws.Cells.Item(48, 3).Value = 59.993
$ws.Cells.Item(49, 1).Value = 45982.48958333334
$ws.Cells.Item(49, 2).Value = 0
$ws.Cells.Item(49, 3).Value = 70.01000000000001
$ws.Cells.Item(50, 1).Value = 45982.5
$ws.Cells.Item(50, 2).Value = 0
$ws.Cells.Item(50, 3).Value = 54.577
$ws.Cells.Item(51, 1).Value = 45982.51041666666
$ws.Cells.Item(51, 2).Value = 4.113
$ws.Cells.Item(51, 3).Value = 9.210000000000001
$ws.Cells.Item(52, 1).Value = 45982.52083333334
$ws.Cells.Item(52, 2).Value = 1.816
$ws.Cells.Item(52, 3).Value = 0.033
$ws.Cells.Item(53, 1).Value = 45982.53125
$ws.Cells.Item(53, 2).Value = 2.059
$ws.Cells.Item(53, 3).Value = 1.585
$ws.Cells.Item(54, 1).Value = 45982.54166666666
$ws.Cells.Item(54, 2).Value = 0.961
$ws.Cells.Item(54, 3).Value = 3.075
$ws.Cells.Item(55, 1).Value = 45982.55208333334
$ws.Cells.Item(55, 2).Value = 15.74
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(56, 1).Value = 45982.5625
$ws.Cells.Item(56, 2).Value = 14.549
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(57, 1).Value = 45982.57291666666
$ws.Cells.Item(57, 2).Value = 4.191
$ws.Cells.Item(57, 3).Value = 2.099
$ws.Cells.Item(58, 1).Value = 45982.58333333334
$ws.Cells.Item(58, 2).Value = 0
$ws.Cells.Item(58, 3).Value = 34.288
$ws.Cells.Item(59, 1).Value = 45982.59375
$ws.Cells.Item(59, 2).Value = 0.473
$ws.Cells.Item(59, 3).Value = 11.651

# New row 60: copy the date format/style from row 59, then set values
$ws.Range("A59").Copy()
$ws.Range("A60").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(60, 1).Value = 45982.60416666666
$ws.Cells.Item(60, 2).Value = 0
$ws.Cells.Item(60, 3).Value = 0
